$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "F4"  = 722
        "F7"  = 2830
        "F9"  = 1826
        "F12" = 748
        "F13" = 899
        "F15" = 378
        "F18" = 51
        "F20" = 6652
        "F21" = 254
        "F22" = 1545
        "F23" = 161
        "F24" = 179
        "F27" = 271
        "F29" = 1101
        "F34" = 780
        "F35" = 1480
        "F39" = 18
        "F41" = 206
    }
    "全部类型" = @{
        "F4"  = 722
        "F10" = 2830
        "F12" = 1826
        "F15" = 748
        "F17" = 899
        "F19" = 378
        "F21" = 51
        "F23" = 6652
        "F24" = 254
        "F25" = 1545
        "F27" = 161
        "F28" = 179
        "F31" = 271
        "F33" = 1101
        "F38" = 780
        "F39" = 1480
        "F43" = 18
        "F45" = 206
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
